# Adds a new "2022-Q4" quarterly sheet to the workbook, right after the
# "总计" (summary) sheet, shifting the existing quarter sheets down by one
# position, and inserts the corresponding new row into the "总计" sheet.

function Set-TextCell($range) {
    # Forces a cell to store its content as literal text even if it looks
    # like a number (mirrors typing `'513980` into Excel), then clears the
    # quote-prefix formatting that the apostrophe trick leaves behind so
    # the cell keeps the workbook's default style.
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert a new row for 2022-Q4 right after the header,
#    shifting every existing data row down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Preserve the row-9 style on the new row 10 before shifting the values.
$summary.Range("A9").Copy($summary.Range("A10"))

for ($r = 9; $r -ge 2; $r--) {
    $src = "B" + $r + ":D" + $r
    $dst = "B" + ($r + 1) + ":D" + ($r + 1)
    $summary.Range($src).Copy($summary.Range($dst))
}

$summary.Range("A10").Value2 = 8
$summary.Range("B2").Value2 = "2022-Q4"
$summary.Range("C2").Value2 = 33
$summary.Range("D2").Value2 = 5.28

# ---------------------------------------------------------------------
# 2. Duplicate the "2022-Q3" sheet (position 2) as a template - it has
#    the right header row/styles - placing the copy right before it, then
#    rename + overwrite its data with the 2022-Q4 fund holdings.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item(2)
$template.Copy($template, $null)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Restore the "last sheet is the active tab" state the copy operation
# disturbed (Excel auto-selects the freshly created sheet).
$wb.Worksheets.Item($wb.Worksheets.Count).Select()

$rows = @(
    @("0", "513980", "景顺长城中证港股通科技ETF", "21.06", "98.84", "6.12", "1.2889", 7),
    @("1", "159636", "工银瑞信国证港股通科技ETF", "11.19", "98.28", "4.44", "0.4968", 6),
    @("2", "513700", "鹏华中证港股通医药卫生综合ETF", "4.34", "95.37", "10.17", "0.4414", 2),
    @("3", "007718", "中银创新医疗混合A", "11.61", "80.72", "3.46", "0.4017", 6),
    @("4", "517120", "华泰柏瑞中证沪港深创新药产业ETF", "3.74", "96.74", "8.09", "0.3026", 4),
    @("5", "513860", "海富通中证港股通科技ETF", "4.55", "96.29", "5.92", "0.2694", 7),
    @("6", "513120", "广发中证香港创新药（QDII-ETF）", "1.85", "98.69", "12.88", "0.2383", 1),
    @("7", "862001", "光大阳光香港精选混合（QDII）A 人民币", "3.90", "92.65", "4.80", "0.1872", 7),
    @("8", "862011", "光大阳光香港精选混合（QDII）A 美元", "3.90", "92.65", "4.80", "0.1872", 7),
    @("9", "862012", "光大阳光香港精选混合（QDII）C 人民币", "3.90", "92.65", "4.80", "0.1872", 7),
    @("10", "513200", "易方达中证港股通医药卫生综合ETF", "1.69", "95.67", "10.46", "0.1768", 2),
    @("11", "513020", "国泰中证港股通科技ETF", "3.00", "94.32", "5.60", "0.1680", 7),
    @("12", "159748", "富国中证沪港深创新药产业ETF", "1.70", "99.36", "8.41", "0.1430", 4),
    @("13", "010500", "中银创新医疗混合C", "3.82", "80.72", "3.46", "0.1322", 6),
    @("14", "517110", "国泰中证沪港深创新药产业ETF", "1.44", "94.70", "7.10", "0.1022", 4),
    @("15", "513150", "华泰柏瑞中证港股通科技ETF", "1.54", "97.23", "6.21", "0.0956", 7),
    @("16", "159776", "银华中证港股通医药卫生综合ETF", "0.81", "93.98", "10.00", "0.0810", 2),
    @("17", "159718", "平安中证港股通医药卫生综合ETF", "0.67", "94.63", "10.04", "0.0673", 2),
    @("18", "159751", "鹏华中证港股通科技ETF", "1.03", "95.29", "5.95", "0.0613", 7),
    @("19", "006787", "泰康中证港股通大消费主题指数C", "1.28", "94.75", "4.09", "0.0524", 7),
    @("20", "011071", "鹏华安悦一年持有期混合A", "5.19", "24.24", "0.79", "0.0410", 4),
    @("21", "014129", "西藏东财中证沪港深创新药产业指数C", "0.33", "95.08", "8.04", "0.0265", 4),
    @("22", "014128", "西藏东财中证沪港深创新药产业指数A", "0.29", "95.08", "8.04", "0.0233", 4),
    @("23", "006786", "泰康中证港股通大消费主题指数A", "0.56", "94.75", "4.09", "0.0229", 7),
    @("24", "012111", "鹏华安颐混合A", "2.28", "25.16", "0.83", "0.0189", 4),
    @("25", "517990", "招商中证沪港深500医药卫生ETF", "0.32", "97.39", "5.37", "0.0172", 5),
    @("26", "005701", "上投摩根香港精选港股通混合A", "0.48", "89.99", "3.24", "0.0156", 4),
    @("27", "860008", "光大阳光生活18个月持有期混合A", "0.33", "90.14", "3.79", "0.0125", 4),
    @("28", "860060", "光大阳光生活18个月持有期混合B", "0.18", "90.14", "3.79", "0.0068", 4),
    @("29", "012112", "鹏华安颐混合C", "0.80", "25.16", "0.83", "0.0066", 4),
    @("30", "860061", "光大阳光生活18个月持有期混合C", "0.08", "90.14", "3.79", "0.0030", 4),
    @("31", "016921", "上投摩根香港精选港股通混合C", "0.02", "89.99", "3.24", "0.0006", 4),
    @("32", "011072", "鹏华安悦一年持有期混合C", "0.05", "24.24", "0.79", "0.0004", 4)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value2 = [int]$row[0]

    $newSheet.Range("B$r").Formula = "'" + $row[1]
    Set-TextCell $newSheet.Range("B$r")

    $newSheet.Range("C$r").Formula = "'" + $row[2]
    Set-TextCell $newSheet.Range("C$r")

    $newSheet.Range("D$r").Formula = "'" + $row[3]
    Set-TextCell $newSheet.Range("D$r")

    $newSheet.Range("E$r").Formula = "'" + $row[4]
    Set-TextCell $newSheet.Range("E$r")

    $newSheet.Range("F$r").Formula = "'" + $row[5]
    Set-TextCell $newSheet.Range("F$r")

    $newSheet.Range("G$r").Formula = "'" + $row[6]
    Set-TextCell $newSheet.Range("G$r")

    $newSheet.Range("H$r").Value2 = [int]$row[7]

    $r = $r + 1
}

# The new sheet now has 33 data rows (2022-Q3 only had 24), clear the
# leftover rows 35-36 if the template had more rows than the new data
# (2022-Q3 template had 25 rows total, new data needs 34, so nothing to
# clear - but guard anyway in case template/content sizes ever diverge).
$lastNewRow = 1 + $rows.Count
$templateLastRow = 25
if ($templateLastRow -gt $lastNewRow) {
    $clearFrom = $lastNewRow + 1
    $newSheet.Range("A" + $clearFrom + ":H" + $templateLastRow).Clear()
}
